# Choosing a solution for predicting fingers
# Insert the "Choose a solution and develop a plan to implement it" heading
# and its explanatory paragraph just before the trailing bookmark/space
# paragraphs at the end of the "Predicting Fingers" section, then fold the
# old trailing blank/space paragraph into the new final paragraph (which
# also keeps carrying the _GoBack bookmark).

$d = $word.ActiveDocument

# Locate the last body paragraph of the "Predicting Fingers" section via its
# unique text, then walk forward paragraph-by-paragraph from there so we
# don't depend on hard-coded paragraph indices.
$anchor = $d.Content
$anchor.Find.Execute(
    "This same solution will work for all cases of determining which finger she will stop on when counting to 10, 100, and 1000.",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchor.Expand(4) | Out-Null               # wdParagraph -> whole paragraph

# The (untouched) blank paragraph right after the anchor paragraph.
$blankPara = $d.Range($anchor.End, $anchor.End)
$blankPara.Expand(4) | Out-Null

# The paragraph that currently holds only the _GoBack bookmark.
$bookmarkPara = $d.Range($blankPara.End, $blankPara.End)
$bookmarkPara.Expand(4) | Out-Null

# The trailing paragraph that currently holds just a single space.
$spacePara = $d.Range($bookmarkPara.End, $bookmarkPara.End)
$spacePara.Expand(4) | Out-Null

$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Choose a solution and develop a plan to implement it</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>The solution will be to count in her manner using my left hand only and when I reach 10, 100, and 1000, the fingers I land on will be the ones she stops at.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

# Replacing the (collapsed) bookmark paragraph in place with the three new
# paragraphs above: a bold heading, the explanatory body paragraph, and a
# final paragraph that keeps the single-space run together with the
# _GoBack bookmark (this is what used to be two separate paragraphs).
$insertionPoint = $d.Range($bookmarkPara.Start, $bookmarkPara.Start)
$insertionPoint.InsertXML($newXml)

# The old trailing "single space" paragraph is now redundant (its content
# now lives in the new final paragraph above), so remove it.
$staleLen = $spacePara.End - $spacePara.Start
$staleSpacePara = $d.Range($insertionPoint.End, $insertionPoint.End)
$staleSpacePara.Expand(4) | Out-Null
$staleSpacePara.Delete()
